# Income statement "yearly" roll-forward:
#  - Drop the oldest fiscal-year column's label/date/values and append a new
#    fiscal-year column (1401/12) with its publish date(s) and figures.
#  - Every data column shifts one position to the left (D<-E<-F<-G<-H) and
#    the freed rightmost column (H) receives the newly reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8: "دوره مالی" period headers (D8:H8) - shift left, new H8 label
# ---------------------------------------------------------------------
$e8 = $ws.Range("E8").Value2
$f8 = $ws.Range("F8").Value2
$g8 = $ws.Range("G8").Value2
$h8 = $ws.Range("H8").Value2

$ws.Range("D8").Value2 = $e8
$ws.Range("E8").Value2 = $f8
$ws.Range("F8").Value2 = $g8
$ws.Range("G8").Value2 = $h8
$ws.Range("H8").Value2 = "12 ماهه منتهی به 1401/12"

# ---------------------------------------------------------------------
# Row 9: "تاریخ انتشار" publish dates (D9:H9). D9:F9 shift left (same as
# every other row); the two newest columns (G9,H9) both get freshly
# reported publish dates for the new filing - the old H9 date is
# superseded entirely, not carried forward.
# ---------------------------------------------------------------------
$e9 = $ws.Range("E9").Value2
$f9 = $ws.Range("F9").Value2
$g9 = $ws.Range("G9").Value2

$ws.Range("D9").Value2 = $e9
$ws.Range("E9").Value2 = $f9
$ws.Range("F9").Value2 = $g9
$ws.Range("G9").Value2 = "1402-02-13 (9)"
$ws.Range("H9").Value2 = "1402-02-13 (2)"

# ---------------------------------------------------------------------
# Data rows 11-27: shift each row's D:H figures left by one column and
# write the freshly reported value into the vacated H column.
# ---------------------------------------------------------------------
$newH = @{
    11 = 5968165
    12 = -3398826
    13 = 2569339
    14 = -256206
    15 = 0
    16 = 0
    17 = 2313133
    18 = -33600
    19 = 255742
    20 = 2535275
    21 = -213046
    22 = 2322229
    23 = 0
    24 = 2322229
    25 = 11611
    26 = 200000
    27 = 11611
}

foreach ($row in 11..27) {
    $e = $ws.Range("E$row").Value2
    $f = $ws.Range("F$row").Value2
    $g = $ws.Range("G$row").Value2
    $h = $ws.Range("H$row").Value2

    $ws.Range("D$row").Value2 = $e
    $ws.Range("E$row").Value2 = $f
    $ws.Range("F$row").Value2 = $g
    $ws.Range("G$row").Value2 = $h
    $ws.Range("H$row").Value2 = $newH[$row]
}
